# Re-simulated Week 17, factoring in more player injuries
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rushing")

# New full data set for the Rushing sheet (rows 2-11); the old sheet had
# rows 2-14 (13 players). Three players (A.Gibson, J.McKissic, L.Thomas)
# were dropped from the simulated results, and the remaining players' stats
# were re-simulated (not just shifted).
$data = @(
    @(0, "T.Heinicke",  10, 10, 13, 9),
    @(1, "K.Allen",      1,  0,  1, 1),
    @(2, "J.Patterson", 135, 85, 31, 46),
    @(3, "W.Smallwood",  24, 18,  9, 6),
    @(4, "J.Williams",   29, 15,  5, 4),
    @(6, "T.McLaurin",    1,  0,  0, 0),
    @(7, "C.Samuel",      3,  0,  1, 0),
    @(8, "D.Brown",       1,  2,  0, 1),
    @(9, "C.Sims",        1,  0,  0, 0),
    @(10, "D.Carter",     6,  3,  0, 1)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $row++
}

# Remove the now-unused trailing rows (old sheet went to row 14).
$lastOldRow = 14
if ($lastOldRow -ge $row) {
    $deleteRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($lastOldRow, 6))
    $deleteRange.Delete()
}

# Match the saved selection state left behind in the authored workbook.
$ws.Activate()
[void]$ws.Range("I10").Select()
